$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cell H1
$ws.Range("H1").Value = "Save"
$ws.Range("H1").Font.Bold = $true
$ws.Range("H1").HorizontalAlignment = -4108  # xlCenter
$ws.Range("H1").VerticalAlignment = -4160    # xlTop
$ws.Range("H1").Borders.LineStyle = 1
$ws.Range("H1").Borders.Weight = 2

# Data values for H2:H8 ("Save" column)
$saveValues = @(1, 0, 0, 1, 1, 0, 0)
for ($i = 0; $i -lt $saveValues.Length; $i++) {
    $row = 2 + $i
    $ws.Cells.Item($row, 8).Value = $saveValues[$i]
}
